# Apply dialog-act re-annotation updates to columns I (DAMSLTag) and J (DialogAct)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 'sd'
$ws.Range("J2").Value = 'Statement-non-opinion'
$ws.Range("I16").Value = 'aa'
$ws.Range("J16").Value = 'Agree/Accept'
$ws.Range("I17").Value = 'sd'
$ws.Range("J17").Value = 'Statement-non-opinion'
$ws.Range("I24").Value = 'aa'
$ws.Range("J24").Value = 'Agree/Accept'
$ws.Range("I25").Value = 'sd'
$ws.Range("J25").Value = 'Statement-non-opinion'
$ws.Range("I37").Value = 'aa'
$ws.Range("J37").Value = 'Agree/Accept'
$ws.Range("I45").Value = 'aa'
$ws.Range("J45").Value = 'Agree/Accept'
$ws.Range("I46").Value = 'ba'
$ws.Range("J46").Value = 'Appreciation'
$ws.Range("I56").Value = 'sd'
$ws.Range("J56").Value = 'Statement-non-opinion'
$ws.Range("I60").Value = 'sv'
$ws.Range("J60").Value = 'Statement-opinion'
$ws.Range("I62").Value = 'sd'
$ws.Range("J62").Value = 'Statement-non-opinion'
$ws.Range("I83").Value = 'ba'
$ws.Range("J83").Value = 'Appreciation'
$ws.Range("I86").Value = 'sv'
$ws.Range("J86").Value = 'Statement-opinion'
$ws.Range("I94").Value = 'aa'
$ws.Range("J94").Value = 'Agree/Accept'
$ws.Range("I97").Value = 'sd'
$ws.Range("J97").Value = 'Statement-non-opinion'
$ws.Range("I99").Value = 'aa'
$ws.Range("J99").Value = 'Agree/Accept'
$ws.Range("I104").Value = 'aa'
$ws.Range("J104").Value = 'Agree/Accept'
$ws.Range("I106").Value = 'sd'
$ws.Range("J106").Value = 'Statement-non-opinion'
$ws.Range("I111").Value = 'sd'
$ws.Range("J111").Value = 'Statement-non-opinion'
$ws.Range("I113").Value = 'aa'
$ws.Range("J113").Value = 'Agree/Accept'
$ws.Range("I119").Value = 'sd'
$ws.Range("J119").Value = 'Statement-non-opinion'
$ws.Range("I127").Value = 'sd'
$ws.Range("J127").Value = 'Statement-non-opinion'
$ws.Range("I133").Value = 'sv'
$ws.Range("J133").Value = 'Statement-opinion'
$ws.Range("I147").Value = 'aa'
$ws.Range("J147").Value = 'Agree/Accept'
$ws.Range("I150").Value = 'aa'
$ws.Range("J150").Value = 'Agree/Accept'
$ws.Range("I153").Value = 'sd'
$ws.Range("J153").Value = 'Statement-non-opinion'
$ws.Range("I155").Value = 'sd'
$ws.Range("J155").Value = 'Statement-non-opinion'
$ws.Range("I158").Value = '%'
$ws.Range("J158").Value = 'Uninterpretable'
$ws.Range("I166").Value = 'sd'
$ws.Range("J166").Value = 'Statement-non-opinion'
$ws.Range("I167").Value = 'sd'
$ws.Range("J167").Value = 'Statement-non-opinion'
$ws.Range("I191").Value = '%'
$ws.Range("J191").Value = 'Uninterpretable'
$ws.Range("I204").Value = 'sd'
$ws.Range("J204").Value = 'Statement-non-opinion'
$ws.Range("I222").Value = 'sv'
$ws.Range("J222").Value = 'Statement-opinion'
$ws.Range("I224").Value = 'aa'
$ws.Range("J224").Value = 'Agree/Accept'
$ws.Range("I235").Value = 'sd'
$ws.Range("J235").Value = 'Statement-non-opinion'
$ws.Range("I239").Value = 'sd'
$ws.Range("J239").Value = 'Statement-non-opinion'
$ws.Range("I242").Value = 'sd'
$ws.Range("J242").Value = 'Statement-non-opinion'
$ws.Range("I248").Value = 'sv'
$ws.Range("J248").Value = 'Statement-opinion'
$ws.Range("I249").Value = 'sd'
$ws.Range("J249").Value = 'Statement-non-opinion'
$ws.Range("I253").Value = 'sd'
$ws.Range("J253").Value = 'Statement-non-opinion'
$ws.Range("I256").Value = 'sd'
$ws.Range("J256").Value = 'Statement-non-opinion'
$ws.Range("I259").Value = 'sd'
$ws.Range("J259").Value = 'Statement-non-opinion'
$ws.Range("I268").Value = 'sv'
$ws.Range("J268").Value = 'Statement-opinion'
$ws.Range("I281").Value = 'sv'
$ws.Range("J281").Value = 'Statement-opinion'
$ws.Range("I283").Value = 'aa'
$ws.Range("J283").Value = 'Agree/Accept'
$ws.Range("I284").Value = 'aa'
$ws.Range("J284").Value = 'Agree/Accept'
$ws.Range("I293").Value = 'sv'
$ws.Range("J293").Value = 'Statement-opinion'
$ws.Range("I297").Value = 'ba'
$ws.Range("J297").Value = 'Appreciation'
$ws.Range("I319").Value = 'ba'
$ws.Range("J319").Value = 'Appreciation'
$ws.Range("I326").Value = 'aa'
$ws.Range("J326").Value = 'Agree/Accept'
$ws.Range("I328").Value = 'sv'
$ws.Range("J328").Value = 'Statement-opinion'
$ws.Range("I331").Value = '%'
$ws.Range("J331").Value = 'Uninterpretable'
$ws.Range("I332").Value = 'b'
$ws.Range("J332").Value = 'Acknowledge (Backchannel)'
$ws.Range("I335").Value = 'sd'
$ws.Range("J335").Value = 'Statement-non-opinion'
$ws.Range("I337").Value = 'sd'
$ws.Range("J337").Value = 'Statement-non-opinion'
$ws.Range("I354").Value = 'sd'
$ws.Range("J354").Value = 'Statement-non-opinion'
$ws.Range("I362").Value = 'sd'
$ws.Range("J362").Value = 'Statement-non-opinion'
$ws.Range("I363").Value = 'sv'
$ws.Range("J363").Value = 'Statement-opinion'
$ws.Range("I364").Value = 'qy'
$ws.Range("J364").Value = 'Yes-No-Question'
$ws.Range("I365").Value = 'qy'
$ws.Range("J365").Value = 'Yes-No-Question'
$ws.Range("I370").Value = 'sd'
$ws.Range("J370").Value = 'Statement-non-opinion'
$ws.Range("I372").Value = 'sd'
$ws.Range("J372").Value = 'Statement-non-opinion'
$ws.Range("I375").Value = 'sd'
$ws.Range("J375").Value = 'Statement-non-opinion'
$ws.Range("I377").Value = 'sv'
$ws.Range("J377").Value = 'Statement-opinion'
$ws.Range("I379").Value = 'sd'
$ws.Range("J379").Value = 'Statement-non-opinion'
$ws.Range("I386").Value = 'sd'
$ws.Range("J386").Value = 'Statement-non-opinion'
$ws.Range("I388").Value = 'sd'
$ws.Range("J388").Value = 'Statement-non-opinion'
$ws.Range("I390").Value = 'sd'
$ws.Range("J390").Value = 'Statement-non-opinion'
$ws.Range("I391").Value = '%'
$ws.Range("J391").Value = 'Uninterpretable'
$ws.Range("I393").Value = 'sv'
$ws.Range("J393").Value = 'Statement-opinion'
$ws.Range("I406").Value = 'sd'
$ws.Range("J406").Value = 'Statement-non-opinion'
$ws.Range("I410").Value = 'sv'
$ws.Range("J410").Value = 'Statement-opinion'
$ws.Range("I411").Value = 'sv'
$ws.Range("J411").Value = 'Statement-opinion'
$ws.Range("I428").Value = 'sd'
$ws.Range("J428").Value = 'Statement-non-opinion'
$ws.Range("I432").Value = 'sd'
$ws.Range("J432").Value = 'Statement-non-opinion'
$ws.Range("I441").Value = 'sd'
$ws.Range("J441").Value = 'Statement-non-opinion'
$ws.Range("I452").Value = 'sd'
$ws.Range("J452").Value = 'Statement-non-opinion'
$ws.Range("I453").Value = 'sv'
$ws.Range("J453").Value = 'Statement-opinion'
$ws.Range("I456").Value = 'sd'
$ws.Range("J456").Value = 'Statement-non-opinion'
$ws.Range("I458").Value = 'sd'
$ws.Range("J458").Value = 'Statement-non-opinion'
$ws.Range("I460").Value = 'sd'
$ws.Range("J460").Value = 'Agree/Accept'
